$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells remain text-typed (matches original inlineStr cells) so
# values like "1.001" or "30.370.59" are not coerced into numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.370.59'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.936.86'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7742'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +7.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '246.39'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.89'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07068'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7827'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08013'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.934.52'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.81'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.359.95'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '255.54'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007982'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.820'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.184.68'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.764'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.569'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.66'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.95%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.08'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.284'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.49%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.519'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.427'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.134'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05173'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.288'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7495'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.774'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01964'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.809'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.80'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.442'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4510'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8362'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.88'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.813'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.517'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '984.98'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +10.78%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4165'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.85%  '
